# Update resultados - 2026-02-14 01:09:33
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column G: "estado" header ---
$ws.Range("G1").Value = "estado"

# --- Row 3-8: shift "jugador" (E) and "total" (F) values up by one row ---
# (the player that used to occupy row N+1 now occupies row N; the old
#  "Luján Martínez, Benjamín" entry in row 3 is dropped)
$ws.Range("E3").Value = "Micheloud, Artemio"
$ws.Range("F3").Value = 20

$ws.Range("E4").Value = "Elcura, Lorenzo"
$ws.Range("F4").Value = 23

$ws.Range("E5").Value = "Bogado, Ogán"
$ws.Range("F5").Value = 33

$ws.Range("E6").Value = "Vera, Bautista"
$ws.Range("F6").Value = 34

$ws.Range("E7").Value = "Esborraz, Juan Cruz"
$ws.Range("F7").Value = 38

$ws.Range("E8").Value = "Núñez, Valentino"
$ws.Range("F8").Value = 41

# --- Row 9: becomes a new "Birdies" / "caballeros" entry with no total (NPT) ---
$ws.Range("B9").Value = "Birdies"
$ws.Range("C9").Value = "caballeros"
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = "Malvasio, Joaquín"
$ws.Range("F9").ClearContents()

# --- Column G "estado" values for every data row ---
$ws.Range("G2").Value = "OK"
$ws.Range("G3").Value = "OK"
$ws.Range("G4").Value = "OK"
$ws.Range("G5").Value = "OK"
$ws.Range("G6").Value = "OK"
$ws.Range("G7").Value = "OK"
$ws.Range("G8").Value = "OK"
$ws.Range("G9").Value = "NPT"
$ws.Range("G10").Value = "OK"
$ws.Range("G11").Value = "OK"
$ws.Range("G12").Value = "OK"
$ws.Range("G13").Value = "OK"
$ws.Range("G14").Value = "OK"
$ws.Range("G15").Value = "OK"
